$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Remove the trailing " 用户" (space + red "用户") runs that follow
#    "用户喜欢什么，就创造什么。" at the end of the first paragraph.
# ---------------------------------------------------------------------
$d.Content.Find.Execute("用户喜欢什么，就创造什么。 用户", $true, $false, $false, $false, $false,
                         $true, 1, $false, "用户喜欢什么，就创造什么。", 2)

# ---------------------------------------------------------------------
# 2) Give the (now last) paragraph's own paragraph-mark run properties
#    the same "hint=eastAsia" that its text run already carries, i.e.
#    turn:
#       <w:pPr><w:rPr><w:rFonts .../></w:rPr></w:pPr>
#    into:
#       <w:pPr><w:rPr><w:rFonts ... w:hint="eastAsia"/></w:rPr></w:pPr>
#    for the "图啥呢，图个快乐。" paragraph. The OM has no direct setter
#    for the pilcrow's own rFonts/hint, so rewrite that paragraph (mark
#    included) via InsertXML with the desired markup, then fold the
#    leftover blank paragraph (InsertXML always leaves one behind when
#    the paragraph mark itself is replaced) back into it.
# ---------------------------------------------------------------------
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$target = $lastPara.Range.Text.TrimEnd([char]13, [char]7)

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq $target) {

        $openXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="18ECADD6" w14:textId="43291BCF" w:rsidR="006759F9" w:rsidRPr="006759F9" w:rsidRDefault="006759F9"><w:pPr><w:rPr><w:rFonts w:asciiTheme="minorEastAsia" w:hAnsiTheme="minorEastAsia" w:cs="Times New Roman" w:hint="eastAsia"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr></w:pPr><w:r w:rsidRPr="006759F9"><w:rPr><w:rFonts w:asciiTheme="minorEastAsia" w:hAnsiTheme="minorEastAsia" w:cs="Times New Roman" w:hint="eastAsia"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>图啥呢，图个快乐。</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

        $beforeCount = $d.Paragraphs.Count
        $p.Range.InsertXML($openXml)

        # InsertXML replaces the paragraph's content/mark but - because the
        # mark itself was part of the replaced range - leaves one extra,
        # now-empty paragraph behind right after it. Delete it by merging
        # it (and the freshly written paragraph's own new mark) together;
        # a collapsed/mark-only range is a no-op in this engine, but a
        # range spanning both marks merges the two paragraphs and keeps
        # the earlier (correct) paragraph's formatting.
        if ($d.Paragraphs.Count -gt $beforeCount) {
            $mergeStart = $p.Range.End - 1
            $mergeRange = $d.Range($mergeStart, $mergeStart + 2)
            $mergeRange.Delete()
        }
        break
    }
}
